{"js": "// Update the worksheet date header and the 25 division-problem answers.\n// Cells/paragraphs are targeted positionally (index-based) rather than by\n// text search, because several of the old problem strings repeat and some\n// new values collide with other cells' old values (e.g. \"46\u00f72=23, 0\"\n// appears twice, and \"41\u00f72=20, 1\" is both an old value and a new value).\n\nconst body = context.document.body;\n\n// --- Date header (first paragraph of the document) -------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2026-01-23 Friday\", Word.InsertLocation.replace);\n\n// --- Table of division problems ---------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New answer for each populated cell, in reading order (row-major,\n// top-to-bottom / left-to-right), matching the data rows 0, 4, 8, 12, 16\n// (0-indexed; the rows in between are blank spacer rows).\nconst answers = [\n  \"30\u00f76=5, 0\",  \"44\u00f77=6, 2\",  \"27\u00f79=3, 0\",  \"66\u00f73=22, 0\", \"10\u00f72=5, 0\",\n  \"27\u00f77=3, 6\",  \"44\u00f74=11, 0\", \"50\u00f78=6, 2\",  \"58\u00f73=19, 1\", \"47\u00f79=5, 2\",\n  \"41\u00f72=20, 1\", \"11\u00f75=2, 1\",  \"15\u00f79=1, 6\",  \"65\u00f79=7, 2\",  \"22\u00f78=2, 6\",\n  \"46\u00f74=11, 2\", \"80\u00f73=26, 2\", \"45\u00f74=11, 1\", \"17\u00f74=4, 1\",  \"23\u00f73=7, 2\",\n  \"35\u00f73=11, 2\", \"20\u00f79=2, 2\",  \"71\u00f72=35, 1\", \"14\u00f74=3, 2\",  \"96\u00f74=24, 0\",\n];\n\nconst dataRows = [0, 4, 8, 12, 16];\nconst cols = [0, 1, 2, 3, 4];\n\nlet idx = 0;\nfor (const row of dataRows) {\n  for (const col of cols) {\n    const cell = table.getCell(row, col);\n    cell.value = answers[idx];\n    idx += 1;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 division-problem answers.\n# Values are targeted positionally (paragraph / table cell) rather than by\n# text search, because several of the old problem strings repeat and some\n# new values collide with other cells' old values (e.g. \"46\u00f72=23, 0\"\n# appears twice, and \"41\u00f72=20, 1\" is both an old value and a new value).\n\n$d = $word.ActiveDocument\n\n# --- Date header (first paragraph of the document) -------------------\n$dateRange = $d.Paragraphs(1).Range\n$dateRange.End = $dateRange.End - 1\n$dateRange.Text = \"2026-01-23 Friday\"\n\n# --- Table of division problems ---------------------------------------\n$t = $d.Tables(1)\n\n# New answer for each populated cell, in reading order (row-major,\n# top-to-bottom / left-to-right), matching the data rows 1, 5, 9, 13, 17.\n$answers = @(\n    \"30\u00f76=5, 0\",   \"44\u00f77=6, 2\",   \"27\u00f79=3, 0\",   \"66\u00f73=22, 0\",  \"10\u00f72=5, 0\",\n    \"27\u00f77=3, 6\",   \"44\u00f74=11, 0\",  \"50\u00f78=6, 2\",   \"58\u00f73=19, 1\",  \"47\u00f79=5, 2\",\n    \"41\u00f72=20, 1\",  \"11\u00f75=2, 1\",   \"15\u00f79=1, 6\",   \"65\u00f79=7, 2\",   \"22\u00f78=2, 6\",\n    \"46\u00f74=11, 2\",  \"80\u00f73=26, 2\",  \"45\u00f74=11, 1\",  \"17\u00f74=4, 1\",   \"23\u00f73=7, 2\",\n    \"35\u00f73=11, 2\",  \"20\u00f79=2, 2\",   \"71\u00f72=35, 1\",  \"14\u00f74=3, 2\",   \"96\u00f74=24, 0\"\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$cols = @(1, 2, 3, 4, 5)\n\n$idx = 0\nforeach ($row in $dataRows) {\n    foreach ($col in $cols) {\n        $cell = $t.Cell($row, $col)\n        $r = $cell.Range\n        # Trim the trailing cell-mark character so we only overwrite the\n        # visible text (and its paragraph mark), keeping run formatting.\n        $r.End = $r.End - 1\n        $r.Text = $answers[$idx]\n        $idx = $idx + 1\n    }\n}\n"}
